$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2412935323383085
$ws.Range("C2").Value = 0.472636815920398
$ws.Range("J2").Value = 0.01243781094527363
$ws.Range("P2").Value = 0.1691542288557214
$ws.Range("S2").Value = 0.1044776119402985
$ws.Range("B3").Value = 0.0202020202020202
$ws.Range("C3").Value = 0.0202020202020202
$ws.Range("J3").Value = 0.0202020202020202
$ws.Range("P3").Value = 0.7777777777777778
$ws.Range("S3").Value = 0.1616161616161616
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.8125
$ws.Range("S4").Value = 0.1458333333333333
$ws.Range("B6").Value = 0.0975609756097561
$ws.Range("D6").Value = 0.006968641114982578
$ws.Range("E6").Value = 0.003484320557491289
$ws.Range("F6").Value = 0.08362369337979095
$ws.Range("J6").Value = 0.156794425087108
$ws.Range("O6").Value = 0.02787456445993031
$ws.Range("Q6").Value = 0.1951219512195122
$ws.Range("R6").Value = 0.05226480836236934
$ws.Range("S6").Value = 0.3763066202090593
$ws.Range("B7").Value = 0.1179245283018868
$ws.Range("D7").Value = 0.01415094339622642
$ws.Range("F7").Value = 0.0330188679245283
$ws.Range("J7").Value = 0.160377358490566
$ws.Range("O7").Value = 0.02830188679245283
$ws.Range("Q7").Value = 0.2358490566037736
$ws.Range("R7").Value = 0.02358490566037736
$ws.Range("S7").Value = 0.3867924528301887
$ws.Range("B8").Value = 0.1026785714285714
$ws.Range("D8").Value = 0.03125
$ws.Range("F8").Value = 0.08928571428571429
$ws.Range("J8").Value = 0.1383928571428572
$ws.Range("O8").Value = 0.02901785714285714
$ws.Range("Q8").Value = 0.1495535714285714
$ws.Range("R8").Value = 0.05803571428571429
$ws.Range("S8").Value = 0.4017857142857143
$ws.Range("B9").Value = 0.1008403361344538
$ws.Range("D9").Value = 0.008403361344537815
$ws.Range("F9").Value = 0.1050420168067227
$ws.Range("J9").Value = 0.1134453781512605
$ws.Range("O9").Value = 0.02941176470588235
$ws.Range("Q9").Value = 0.1848739495798319
$ws.Range("R9").Value = 0.07983193277310924
$ws.Range("S9").Value = 0.3781512605042017
$ws.Range("B10").Value = 0.1186556927297668
$ws.Range("D10").Value = 0.01989026063100137
$ws.Range("F10").Value = 0.07818930041152264
$ws.Range("J10").Value = 0.1262002743484225
$ws.Range("O10").Value = 0.02606310013717421
$ws.Range("Q10").Value = 0.2263374485596708
$ws.Range("R10").Value = 0.05075445816186557
$ws.Range("S10").Value = 0.3539094650205761
$ws.Range("G11").Value = 0.1671469740634006
$ws.Range("J11").Value = 0.09798270893371758
$ws.Range("K11").Value = 0.2161383285302594
$ws.Range("L11").Value = 0.5014409221902018
$ws.Range("S11").Value = 0.01729106628242075
$ws.Range("G12").Value = 0.7784090909090909
$ws.Range("J12").Value = 0.1534090909090909
$ws.Range("K12").Value = 0.01136363636363636
$ws.Range("L12").Value = 0.02272727272727273
$ws.Range("S12").Value = 0.03409090909090909
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.03308823529411765
$ws.Range("H15").Value = 0.1323529411764706
$ws.Range("I15").Value = 0.08455882352941177
$ws.Range("J15").Value = 0.3860294117647059
$ws.Range("K15").Value = 0.0661764705882353
$ws.Range("M15").Value = 0.01470588235294118
$ws.Range("O15").Value = 0.04779411764705882
$ws.Range("S15").Value = 0.2352941176470588
$ws.Range("F16").Value = 0.02766798418972332
$ws.Range("H16").Value = 0.1541501976284585
$ws.Range("I16").Value = 0.08695652173913043
$ws.Range("J16").Value = 0.4545454545454545
$ws.Range("K16").Value = 0.1264822134387352
$ws.Range("M16").Value = 0.003952569169960474
$ws.Range("O16").Value = 0.05928853754940711
$ws.Range("S16").Value = 0.08695652173913043
$ws.Range("F17").Value = 0.01841620626151013
$ws.Range("H17").Value = 0.1546961325966851
$ws.Range("I17").Value = 0.08839779005524862
$ws.Range("J17").Value = 0.4604051565377532
$ws.Range("K17").Value = 0.08839779005524862
$ws.Range("M17").Value = 0.009208103130755065
$ws.Range("N17").Value = 0.001841620626151013
$ws.Range("O17").Value = 0.06261510128913444
$ws.Range("S17").Value = 0.1160220994475138
$ws.Range("F18").Value = 0.0364963503649635
$ws.Range("H18").Value = 0.1678832116788321
$ws.Range("I18").Value = 0.1532846715328467
$ws.Range("J18").Value = 0.4160583941605839
$ws.Range("K18").Value = 0.05109489051094891
$ws.Range("O18").Value = 0.06569343065693431
$ws.Range("S18").Value = 0.1094890510948905
$ws.Range("F19").Value = 0.01398086828550405
$ws.Range("H19").Value = 0.1979396615158205
$ws.Range("I19").Value = 0.09050772626931568
$ws.Range("J19").Value = 0.3811626195732156
$ws.Range("K19").Value = 0.1140544518027962
$ws.Range("M19").Value = 0.02060338484179544
$ws.Range("N19").Value = 0.0007358351729212656
$ws.Range("O19").Value = 0.07579102281089035
$ws.Range("S19").Value = 0.105224429727741
